# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N - this shifts the old N/O/P "Original/Outstanding" columns one
# place to the right (N->O, O->P, P->Q) and widens the dimension to A1:Q15.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N:N").Insert()

# The newly inserted column inherits the width of the column to its left
# (column M, "Outstanding").
$ws.Columns(14).ColumnWidth = $ws.Columns(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab with S8 selected (this
# also clears the tabSelected flag that used to be on "Acc_Repayment1").
$ws.Activate()
$ws.Range("S8").Select()
